$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("system_chatbot_UI_text")

# Add a new FAQ row for "get similar when not found" error message
# Cell entry order matches the shared-string insertion order recorded in the
# original commit: text_key, then description, then text_value.
$ws.Range("A12").Value = 10
$ws.Range("E12").Value = "error_message"
$ws.Range("G12").Value = "Error message shown when there is no relevant data for the user's question. Encourages the user to ask a more specific question."
$ws.Range("F12").Value = "Sorry, there is no data available related to your question. If you could ask a more specific question, I believe I can assist you better."
